# Complete update of the estimates
#
# The old "HU_E2a" sheet (small 8x8 estimates table) is dropped, and the
# sheet that was "HU_E2a (2)" (the full/updated 15x16 estimates table)
# takes over the name "HU_E2a" and becomes the selected/active sheet.

$wb = $excel.ActiveWorkbook

# Drop the outdated HU_E2a estimates sheet.
$oldEstimates = $wb.Worksheets.Item("HU_E2a")
$oldEstimates.Delete()

# The completed estimates sheet ("HU_E2a (2)") becomes the one and only
# HU_E2a sheet going forward.
$newEstimates = $wb.Worksheets.Item("HU_E2a (2)")
$newEstimates.Name = "HU_E2a"
